$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename sheets ---
$ws1.Name = "liver.nii.gz"
$ws2.Name = "tumor.nii.gz"

# --- Update cell values on sheet1 (liver.nii.gz) ---
$ws1.Range("B2").Value = 0.32957
$ws1.Range("C2").Value = 0.84456
$ws1.Range("D2").Value = 0.89177
$ws1.Range("B3").Value = 0.02628
$ws1.Range("C3").Value = 0.90146
$ws1.Range("D3").Value = 0.88372
$ws1.Range("B4").Value = 0.20106
$ws1.Range("C4").Value = 0.90036
$ws1.Range("D4").Value = 0.9169
$ws1.Range("B5").Value = 0.32613
$ws1.Range("C5").Value = 0.94055
$ws1.Range("D5").Value = 0.93485
$ws1.Range("C6").Value = 0.88209
$ws1.Range("D6").Value = 0.93731
$ws1.Range("B7").Value = 0.0178
$ws1.Range("C7").Value = 0.87438
$ws1.Range("D7").Value = 0.85081
$ws1.Range("B8").Value = 0.3379
$ws1.Range("C8").Value = 0.80693
$ws1.Range("D8").Value = 0.89898
$ws1.Range("B9").Value = 0.16143
$ws1.Range("C9").Value = 0.87774
$ws1.Range("D9").Value = 0.87326
$ws1.Range("B10").Value = 0.38966
$ws1.Range("C10").Value = 0.8646
$ws1.Range("D10").Value = 0.87356
$ws1.Range("B11").Value = 0.0982
$ws1.Range("C11").Value = 0.80365
$ws1.Range("D11").Value = 0.90285
$ws1.Range("C12").Value = 0.80365
$ws1.Range("D12").Value = 0.85081
$ws1.Range("B13").Value = 0.38966
$ws1.Range("C13").Value = 0.94055
$ws1.Range("D13").Value = 0.93731
$ws1.Range("B14").Value = 0.1898075
$ws1.Range("C14").Value = 0.8700433333333333
$ws1.Range("D14").Value = 0.8960108333333333
$ws1.Range("B15").Value = 0.1898075
$ws1.Range("C15").Value = 0.87438
$ws1.Range("D15").Value = 0.8960108333333333

# --- Update cell values on sheet2 (tumor.nii.gz) ---
$ws2.Range("B2").Value = 0.13816
$ws2.Range("C2").Value = 0.78266
$ws2.Range("D2").Value = 0.87473
$ws2.Range("C3").Value = 0.90982
$ws2.Range("D3").Value = 0.89295
$ws2.Range("C4").Value = 0.80014
$ws2.Range("D4").Value = 0.74444
$ws2.Range("B5").Value = 0.24222
$ws2.Range("C5").Value = 0.95806
$ws2.Range("D5").Value = 0.93389
$ws2.Range("C6").Value = 0.83334
$ws2.Range("D6").Value = 0.92549
$ws2.Range("C7").Value = 0.56569
$ws2.Range("D7").Value = 0.04987
$ws2.Range("C8").Value = 0.67168
$ws2.Range("D8").Value = 0.67521
$ws2.Range("B9").Value = 0.10535
$ws2.Range("C9").Value = 0.80164
$ws2.Range("D9").Value = 0.79902
$ws2.Range("B10").Value = 0.25815
$ws2.Range("C10").Value = 0.80407
$ws2.Range("D10").Value = 0.82361
$ws2.Range("B11").Value = 0
$ws2.Range("C11").Value = 0.6596
$ws2.Range("D11").Value = 0.65358
$ws2.Range("C12").Value = 0.56569
$ws2.Range("D12").Value = 0.04987
$ws2.Range("B13").Value = 0.25815
$ws2.Range("C13").Value = 0.95806
$ws2.Range("D13").Value = 0.93389
$ws2.Range("B14").Value = 0.0835025
$ws2.Range("C14").Value = 0.7758708333333333
$ws2.Range("D14").Value = 0.6963791666666667
$ws2.Range("B15").Value = 0
$ws2.Range("C15").Value = 0.80014
$ws2.Range("D15").Value = 0.79902

# --- Column widths: split col B off from col C (bestFit shrank col B) ---
$ws1.Columns.Item(2).ColumnWidth = 9.5
$ws2.Columns.Item(2).ColumnWidth = 9.5

# --- Selection / active sheet state ---
$ws2.Range("F16").Select() | Out-Null
$ws1.Activate()
$ws1.Range("D20").Select() | Out-Null
$ws2.Activate()
$ws2.Range("H12").Select() | Out-Null
$ws1.Activate()
